# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh described in the commit diff.
# For each changed cell (columns H-N: price/profit calculations refreshed from
# current market-board data), set the new numeric value via Range.Value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4430.1875
$ws.Range("I40").Value = 2075
$ws.Range("J40").Value = 4766.643
$ws.Range("K40").Value = 2075
$ws.Range("L40").Value = 4766.643
$ws.Range("M40").Value = -1900
$ws.Range("N40").Value = -5116.643
$ws.Range("H47").Value = 44074
$ws.Range("J47").Value = 44074
$ws.Range("L47").Value = 44074
$ws.Range("N47").Value = -46018
$ws.Range("H54").Value = 21153
$ws.Range("J54").Value = 20084
$ws.Range("L54").Value = 20084
$ws.Range("N54").Value = -21056
$ws.Range("H107").Value = 1600.579
$ws.Range("I107").Value = 1091.381
$ws.Range("K107").Value = 1091.381
$ws.Range("M107").Value = 828.6189999999999
$ws.Range("H111").Value = 16450
$ws.Range("J111").Value = 3281.3333
$ws.Range("L111").Value = 9843.999899999999
$ws.Range("N111").Value = -15977.9999
$ws.Range("H137").Value = 2373.244
$ws.Range("I137").Value = 2711.3914
$ws.Range("J137").Value = 1941.1666
$ws.Range("K137").Value = 8134.174199999999
$ws.Range("L137").Value = 5823.4998
$ws.Range("M137").Value = -5584.174199999999
$ws.Range("N137").Value = -10923.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6367.8115
$ws.Range("I32").Value = 6205.524
$ws.Range("K32").Value = 6205.524
$ws.Range("M32").Value = -5918.524
$ws.Range("H61").Value = 6123
$ws.Range("I61").Value = 4196.019
$ws.Range("J61").Value = 10677.682
$ws.Range("K61").Value = 4196.019
$ws.Range("L61").Value = 10677.682
$ws.Range("M61").Value = -3984.019
$ws.Range("N61").Value = -11101.682
$ws.Range("H110").Value = 1031.4286
$ws.Range("I110").Value = 704.2
$ws.Range("J110").Value = 1849.5
$ws.Range("K110").Value = 704.2
$ws.Range("L110").Value = 1849.5
$ws.Range("M110").Value = 1340.8
$ws.Range("N110").Value = -5939.5
$ws.Range("H132").Value = 1457.8306
$ws.Range("I132").Value = 1275.0209
$ws.Range("K132").Value = 3825.0627
$ws.Range("M132").Value = -1295.0627
$ws.Range("H136").Value = 6123
$ws.Range("I136").Value = 4196.019
$ws.Range("J136").Value = 10677.682
$ws.Range("K136").Value = 12588.057
$ws.Range("L136").Value = 32033.046
$ws.Range("M136").Value = -10038.057
$ws.Range("N136").Value = -37133.046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 401.4
$ws.Range("I22").Value = 389.25
$ws.Range("K22").Value = 389.25
$ws.Range("M22").Value = -216.25
$ws.Range("H99").Value = 7501.75
$ws.Range("J99").Value = 8336
$ws.Range("L99").Value = 8336
$ws.Range("N99").Value = -11332
$ws.Range("H105").Value = 5400.0625
$ws.Range("J105").Value = 8716.166999999999
$ws.Range("L105").Value = 8716.166999999999
$ws.Range("N105").Value = -12210.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2733.932
$ws.Range("I31").Value = 1959.7142
$ws.Range("J31").Value = 4088.8125
$ws.Range("K31").Value = 1959.7142
$ws.Range("L31").Value = 4088.8125
$ws.Range("M31").Value = -1664.7142
$ws.Range("N31").Value = -4678.8125
$ws.Range("H34").Value = 2733.932
$ws.Range("I34").Value = 1959.7142
$ws.Range("J34").Value = 4088.8125
$ws.Range("K34").Value = 1959.7142
$ws.Range("L34").Value = 4088.8125
$ws.Range("M34").Value = -1757.7142
$ws.Range("N34").Value = -4492.8125
$ws.Range("H122").Value = 4182.28
$ws.Range("J122").Value = 3721.1428
$ws.Range("L122").Value = 11163.4284
$ws.Range("N122").Value = -16063.4284
$ws.Range("H125").Value = 90000
$ws.Range("J125").Value = 90000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -94920
$ws.Range("H132").Value = 1336626.9
$ws.Range("H134").Value = 3473.45
$ws.Range("I134").Value = 946.2162
$ws.Range("J134").Value = 7539
$ws.Range("K134").Value = 2838.6486
$ws.Range("L134").Value = 22617
$ws.Range("M134").Value = -303.6486
$ws.Range("N134").Value = -27687
$ws.Range("H139").Value = 18838.625
$ws.Range("I139").Value = 18838.625
$ws.Range("K139").Value = 18838.625
$ws.Range("M139").Value = -13698.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 8573032
$ws.Range("I32").Value = 20000416
$ws.Range("J32").Value = 2494.5
$ws.Range("K32").Value = 60001248
$ws.Range("L32").Value = 7483.5
$ws.Range("M32").Value = -60000965
$ws.Range("N32").Value = -8049.5
$ws.Range("H51").Value = 4136
$ws.Range("I51").Value = 4136
$ws.Range("K51").Value = 12408
$ws.Range("M51").Value = -11948
$ws.Range("H99").Value = 1300
$ws.Range("I99").Value = 1600
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 4800
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -2554
$ws.Range("N99").Value = -7492
$ws.Range("H113").Value = 1576.3572
$ws.Range("I113").Value = 611
$ws.Range("J113").Value = 1692.2
$ws.Range("K113").Value = 1833
$ws.Range("L113").Value = 5076.6
$ws.Range("M113").Value = 337
$ws.Range("N113").Value = -9416.6
$ws.Range("H119").Value = 6663.3335
$ws.Range("I119").Value = 6663.3335
$ws.Range("K119").Value = 19990.0005
$ws.Range("M119").Value = -15152.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10218.737
$ws.Range("I46").Value = 5887.75
$ws.Range("K46").Value = 5887.75
$ws.Range("M46").Value = -5699.75
$ws.Range("H55").Value = 230.1
$ws.Range("I55").Value = 246.77777
$ws.Range("K55").Value = 246.77777
$ws.Range("M55").Value = -73.77777
$ws.Range("H100").Value = 3685.3125
$ws.Range("I100").Value = 3269.5454
$ws.Range("J100").Value = 4600
$ws.Range("K100").Value = 3269.5454
$ws.Range("L100").Value = 4600
$ws.Range("M100").Value = -2728.5454
$ws.Range("N100").Value = -5682
$ws.Range("H122").Value = 4558.952
$ws.Range("I122").Value = 3414.5881
$ws.Range("J122").Value = 9422.5
$ws.Range("K122").Value = 10243.7643
$ws.Range("L122").Value = 28267.5
$ws.Range("M122").Value = -7793.764299999999
$ws.Range("N122").Value = -33167.5
$ws.Range("H132").Value = 3285.169
$ws.Range("I132").Value = 3273.578
$ws.Range("J132").Value = 3305.2307
$ws.Range("K132").Value = 9820.734
$ws.Range("L132").Value = 9915.6921
$ws.Range("M132").Value = -7290.734
$ws.Range("N132").Value = -14975.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 312
$ws.Range("I14").Value = 300
$ws.Range("J14").Value = 336
$ws.Range("K14").Value = 300
$ws.Range("L14").Value = 336
$ws.Range("M14").Value = -132
$ws.Range("N14").Value = -672
$ws.Range("H96").Value = 3800.5
$ws.Range("I96").Value = 3001
$ws.Range("K96").Value = 3001
$ws.Range("M96").Value = -1628
$ws.Range("H100").Value = 969.86664
$ws.Range("I100").Value = 899.8889
$ws.Range("J100").Value = 1599.6666
$ws.Range("K100").Value = 1799.7778
$ws.Range("L100").Value = 3199.3332
$ws.Range("M100").Value = -1258.7778
$ws.Range("N100").Value = -4281.3332
$ws.Range("H107").Value = 1752.2
$ws.Range("I107").Value = 1383.6666
$ws.Range("J107").Value = 2305
$ws.Range("K107").Value = 4150.9998
$ws.Range("L107").Value = 6915
$ws.Range("M107").Value = -2230.9998
$ws.Range("N107").Value = -10755
$ws.Range("H113").Value = 5213936.5
$ws.Range("I113").Value = 27801380
$ws.Range("J113").Value = 1449.9231
$ws.Range("K113").Value = 83404140
$ws.Range("L113").Value = 4349.7693
$ws.Range("M113").Value = -83401970
$ws.Range("N113").Value = -8689.7693
$ws.Range("H132").Value = 1133.2609
$ws.Range("I132").Value = 970.02563
$ws.Range("K132").Value = 2910.07689
$ws.Range("M132").Value = -380.0768899999998
